# Updated cryptos list on Sat Apr 29 11:04:45 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns with the latest scrape,
# and reflects TRON overtaking BinanceUSD in the ranking (rows 16/17 swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the data range to Text format so that values such as
# "1.007" are written back as strings (matching the original inline-string
# cells) instead of being auto-converted into numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.439.37"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.908.30"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "325.08"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "0.4816"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("D8").Value = "0.4062"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.08149"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").Value = "23.44"
$ws.Range("E11").Value = "  +3.09%  "
$ws.Range("D12").Value = "1.927.91"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "6.003"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Value = "7.158"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "90.25"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D18").Value = "0.00001032"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "17.67"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "29.468.75"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "5.621"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").Value = "2.176"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "2.131.08"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").Value = "156.14"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "6.454"
$ws.Range("E27").Value = "  +7.18%  "
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").Value = "2.108"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "120.15"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "1.024"
$ws.Range("E31").Value = "  -4.55%  "
$ws.Range("D32").Value = "0.09518"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "5.508"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").Value = "3.559"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "1.389"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").Value = "0.02269"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "0.06103"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "1.173"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "10.81"
$ws.Range("E39").Value = "  +7.00%  "
$ws.Range("D40").Value = "0.5945"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "7.976"
$ws.Range("E41").Value = "  -3.13%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "1.274"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "2.380"
$ws.Range("E44").Value = "  -4.84%  "
$ws.Range("D45").Value = "12.50"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").Value = "0.07603"
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("D47").Value = "0.5563"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "1.945"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").Value = "116.40"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("D50").Value = "72.54"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").Value = "2.406"
$ws.Range("E51").Value = "  +2.38%  "

# Rows 16 and 17: TRON and BinanceUSD swap rank positions with updated values
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.06789"
$ws.Range("E16").Value = "  +2.50%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.61%  "

# Restore the original (default) cell style now that the text values are set.
$dataRange.Style = "Normal"
